# Generate Report for Handoff
# - Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps for the rows that were just handed off.
# - Marks the Priority column as "ht" (handoff type) for those same rows.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows (on each localization sheet / the matching Overview rows) that are
# part of this handoff batch.
$rows = @(7, 8, 9, 11, 12, 14)

$newHoGenerateDate = "2016-08-18 06:21:03"
$newHandoffDate    = "2016-08-18 06:20:56"

foreach ($r in $rows) {
    # Overview sheet: "Latest HO Xliff Generate Date" (column G)
    $overview.Range("G$r").Value = $newHoGenerateDate

    # de-de sheet: "Latest Handoff Datetime" (column H) mirrors the
    # Overview generate date for this batch.
    $dede.Range("H$r").Value = $newHoGenerateDate

    # zh-cn sheet: "Latest Handoff Datetime" (column H)
    $zhcn.Range("H$r").Value = $newHandoffDate

    # Priority (column E) flips from blank to "ht" on both locale sheets.
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}
